# Auto-generated edit script: apply 2023-03-21 daily crime-count update
# across "Citywide Totals", "By Neighborhood", and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 1391
$ws.Range('J3').Value = 1469
$ws.Range('J4').Value = 330
$ws.Range('J5').Value = 105
$ws.Range('I6').Value = 8969
$ws.Range('J6').Value = 1926
$ws.Range('I7').Value = 26201
$ws.Range('J7').Value = 5221

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('J6').Value = 6
$ws.Range('J7').Value = 19

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 48
$ws.Range('J7').Value = 179

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 23
$ws.Range('J6').Value = 19
$ws.Range('J7').Value = 64

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 78
$ws.Range('J4').Value = 10
$ws.Range('J7').Value = 186

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J3').Value = 30
$ws.Range('J6').Value = 55
$ws.Range('J7').Value = 132

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J6').Value = 49
$ws.Range('J7').Value = 148
$ws.Range('J8').Value = 320
$ws.Range('J9').Value = 33
$ws.Range('J13').Value = 9
$ws.Range('J15').Value = 63
$ws.Range('J19').Value = 178
$ws.Range('J21').Value = 11
$ws.Range('J23').Value = 48
$ws.Range('I29').Value = 1554
$ws.Range('J30').Value = 19
$ws.Range('J33').Value = 224
$ws.Range('J36').Value = 79
$ws.Range('J37').Value = 179
$ws.Range('J42').Value = 210
$ws.Range('J44').Value = 42
$ws.Range('J52').Value = 113
$ws.Range('J53').Value = 49
$ws.Range('J54').Value = 105
$ws.Range('J57').Value = 21
$ws.Range('J63').Value = 23
$ws.Range('J65').Value = 132
$ws.Range('J67').Value = 186
$ws.Range('J71').Value = 23
$ws.Range('J73').Value = 46
$ws.Range('J78').Value = 70
$ws.Range('J79').Value = 162
$ws.Range('J83').Value = 127
$ws.Range('J85').Value = 237
$ws.Range('J90').Value = 61
$ws.Range('J97').Value = 32
$ws.Range('J99').Value = 64
$ws.Range('I101').Value = 26201
$ws.Range('J101').Value = 5221

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J2').Value = 39
$ws.Range('J3').Value = 41
$ws.Range('J7').Value = 127

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 56
$ws.Range('J3').Value = 63
$ws.Range('J4').Value = 10
$ws.Range('J6').Value = 88
$ws.Range('J7').Value = 224

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J6').Value = 51
$ws.Range('J7').Value = 105

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J4').Value = 12
$ws.Range('I6').Value = 434
$ws.Range('J6').Value = 76
$ws.Range('I7').Value = 1554

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 42
$ws.Range('J7').Value = 178

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J2').Value = 15
$ws.Range('J7').Value = 42

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J3').Value = 93
$ws.Range('J4').Value = 15
$ws.Range('J6').Value = 66
$ws.Range('J7').Value = 237

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J6').Value = 17
$ws.Range('J7').Value = 49

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J2').Value = 43
$ws.Range('J3').Value = 39
$ws.Range('J4').Value = 10
$ws.Range('J6').Value = 115
$ws.Range('J7').Value = 210

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('J5').Value = 5
$ws.Range('J6').Value = 9

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J3').Value = 23
$ws.Range('J5').Value = 2
$ws.Range('J7').Value = 70

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('J2').Value = 17
$ws.Range('J3').Value = 14
$ws.Range('J4').Value = 4
$ws.Range('J7').Value = 48

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('J2').Value = 2
$ws.Range('J7').Value = 11

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J3').Value = 60
$ws.Range('J4').Value = 11
$ws.Range('J6').Value = 45
$ws.Range('J7').Value = 162

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J2').Value = 26
$ws.Range('J7').Value = 79

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J6').Value = 45
$ws.Range('J7').Value = 113

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J3').Value = 16
$ws.Range('J7').Value = 63

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('J6').Value = 15
$ws.Range('J7').Value = 33

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J3').Value = 16
$ws.Range('J7').Value = 46

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J3').Value = 3
$ws.Range('J7').Value = 32

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 105
$ws.Range('J6').Value = 87
$ws.Range('J7').Value = 320

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J2').Value = 5
$ws.Range('J4').Value = 3

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J2').Value = 17
$ws.Range('J7').Value = 61

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('J6').Value = 12
$ws.Range('J7').Value = 21

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J2').Value = 8
$ws.Range('J7').Value = 49

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('J3').Value = 4
$ws.Range('J7').Value = 23

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 52
$ws.Range('J7').Value = 148

Write-Output "Applied 126 cell updates."